# Hybrid algorithm (Genetic + Hill Climbing) results for Rosenbrock and SHCB functions
# Updates the second results block (rows 79-108): Griewangk col G, Rosenbrock cols D/J/K
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Griewangk column G (rows 79-108): new values + scientific-notation format ---
$ws.Range("G79").Value = 0.000507244155065755
$ws.Range("G80").Value = 0.000492042623911004
$ws.Range("G81").Value = 0.000238231036261638
$ws.Range("G82").Value = 0.00074985130940941
$ws.Range("G83").Value = 0.000605718384853903
$ws.Range("G84").Value = 0.000165031461983544
$ws.Range("G85").Value = 0.000258484055695062
$ws.Range("G86").Value = 0.000359609351310297
$ws.Range("G87").Value = 0.000212690094702594
$ws.Range("G88").Value = 0.000318842309888323
$ws.Range("G89").Value = 0.000314353052500382
$ws.Range("G90").Value = 0.000590351606529204
$ws.Range("G91").Value = 0.000216414811278431
$ws.Range("G92").Value = 0.000406162110983432
$ws.Range("G93").Value = 0.000338797812932223
$ws.Range("G94").Value = 0.00051788309756573
$ws.Range("G95").Value = 0.000410052814724393
$ws.Range("G96").Value = 0.000464818711093384
$ws.Range("G97").Value = 0.000268561641261655
$ws.Range("G98").Value = 0.000777369128636018
$ws.Range("G99").Value = 0.000573595979664776
$ws.Range("G100").Value = 0.00026446997483609
$ws.Range("G101").Value = 0.000528655358776219
$ws.Range("G102").Value = 0.00045604072869676
$ws.Range("G103").Value = 0.000334071502039345
$ws.Range("G104").Value = 0.000451385928904701
$ws.Range("G105").Value = 0.000434092552727194
$ws.Range("G106").Value = 0.000122958314796584
$ws.Range("G107").Value = 0.000328319407339794
$ws.Range("G108").Value = 0.000360907261502063
$ws.Range("G79:G108").NumberFormat = "0.00E+00"

# --- Rosenbrock column J (rows 79-108) ---
$ws.Range("J79").Value = 0.509152908042577
$ws.Range("J80").Value = 0.00626133650722246
$ws.Range("J81").Value = 0.00626133650722246
$ws.Range("J82").Value = 0.0841707759837449
$ws.Range("J83").Value = 0.00626133650722246
$ws.Range("J84").Value = 0.0841707759837449
$ws.Range("J85").Value = 0.00626133650722246
$ws.Range("J86").Value = 0.0841707759837449
$ws.Range("J87").Value = 0.566539869671253
$ws.Range("J88").Value = 0.00626133650722246
$ws.Range("J89").Value = 0.00626133650722246
$ws.Range("J90").Value = 0.00626133650722246
$ws.Range("J91").Value = 0.00626133650722246
$ws.Range("J92").Value = 0.00626133650722246
$ws.Range("J93").Value = 0.0841707759837449
$ws.Range("J94").Value = 0.00626133650722246
$ws.Range("J95").Value = 0.00626133650722246
$ws.Range("J96").Value = 0.00626133650722246
$ws.Range("J97").Value = 0.00626133650722246
$ws.Range("J98").Value = 0.00626133650722246
$ws.Range("J99").Value = 0.00626133650722246
$ws.Range("J100").Value = 0.00626133650722246
$ws.Range("J101").Value = 0.00626133650722246
$ws.Range("J102").Value = 0.00626133650722246
$ws.Range("J103").Value = 0.124406854575262
$ws.Range("J104").Value = 0.00626133650722246
$ws.Range("J105").Value = 0.0841707759837449
$ws.Range("J106").Value = 0.00626133650722246
$ws.Range("J107").Value = 0.00626133650722246
$ws.Range("J108").Value = 0.00626133650722246

# --- Rosenbrock column K (rows 79-108) ---
$ws.Range("K79").Value = 1.69148420809235
$ws.Range("K80").Value = 4.50248640870764
$ws.Range("K81").Value = 6.30019220924485
$ws.Range("K82").Value = 1.55227442443659
$ws.Range("K83").Value = 1.16001548653857
$ws.Range("K84").Value = 3.69685467612703
$ws.Range("K85").Value = 1.39518541616775
$ws.Range("K86").Value = 0.164116882413045
$ws.Range("K87").Value = 1.47436498496007
$ws.Range("K88").Value = 6.8006306533828
$ws.Range("K89").Value = 5.72760628554448
$ws.Range("K90").Value = 1.27385485712507
$ws.Range("K91").Value = 1.11358080788158
$ws.Range("K92").Value = 0.930695547093326
$ws.Range("K93").Value = 4.36294736351554
$ws.Range("K94").Value = 0.943580494710764
$ws.Range("K95").Value = 0.924905543412075
$ws.Range("K96").Value = 6.30019220924485
$ws.Range("K97").Value = 1.10138273595738
$ws.Range("K98").Value = 5.72163845852264
$ws.Range("K99").Value = 5.7653149210358
$ws.Range("K100").Value = 7.15407970637137
$ws.Range("K101").Value = 4.04955026143054
$ws.Range("K102").Value = 5.85934245279365
$ws.Range("K103").Value = 1.12008916381485
$ws.Range("K104").Value = 3.56140341865566
$ws.Range("K105").Value = 1.33515520130431
$ws.Range("K106").Value = 0.652275979781803
$ws.Range("K107").Value = 4.72751770591243
$ws.Range("K108").Value = 6.67402340879179

# --- Rosenbrock column D (rows 91-108, stdev values recomputed) ---
$ws.Range("D91").Value = 30.0432171759642
$ws.Range("D92").Value = 30.1918842256819
$ws.Range("D93").Value = 24.340199721112
$ws.Range("D94").Value = 29.602187068574
$ws.Range("D95").Value = 31.3351870710883
$ws.Range("D96").Value = 35.6787276554778
$ws.Range("D97").Value = 27.9843126231435
$ws.Range("D98").Value = 27.9087809379636
$ws.Range("D99").Value = 33.1237190184901
$ws.Range("D100").Value = 22.3190475417941
$ws.Range("D101").Value = 29.5337408310067
$ws.Range("D102").Value = 26.4039444560292
$ws.Range("D103").Value = 23.716101736862
$ws.Range("D104").Value = 15.1120804571227
$ws.Range("D105").Value = 25.6577933023307
$ws.Range("D106").Value = 21.1765871401658
$ws.Range("D107").Value = 15.3523044275243
$ws.Range("D108").Value = 29.777237208498

# --- Restore active selection to match the saved view state ---
$ws.Range("D111").Select()
